$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 'BTC', 'Bitcoin', 43078, 844574158403, 8579251721, 0.10174),
    @(2, 'ETH', 'Ethereum', 2309.29, 277313827745, 4495841424, 0.17041),
    @(3, 'USDT', 'Tether', 1, 96197324041, 15934032561, 0.07011000000000001),
    @(4, 'BNB', 'BNB', 300.14, 46154995940, 298877552, -0.5301900000000001),
    @(5, 'SOL', 'Solana', 98.25, 42807942713, 1521841560, -0.5954700000000001),
    @(6, 'XRP', 'XRP', 0.51273, 27905941067, 587870855, -2.02791),
    @(7, 'USDC', 'USDC', 1.001, 27077350476, 2059829225, 0.07106),
    @(8, 'STETH', 'Lido Staked Ether', 2307.44, 21891609938, 8995410, 0.10357),
    @(9, 'ADA', 'Cardano', 0.508147, 17822147261, 234331370, -2.82757),
    @(10, 'AVAX', 'Avalanche', 36.09, 13227179549, 401585307, 0.97945),
    @(11, 'DOGE', 'Dogecoin', 0.079114, 11310731865, 182765076, 0.03996),
    @(13, 'LINK', 'Chainlink', 18.16, 10652043548, 649591385, 1.75003),
    @(12, 'TRX', 'TRON', 0.118736, 10457767911, 185906350, 1.97329),
    @(14, 'DOT', 'Polkadot', 6.81, 9031877303, 95513990, -1.46021),
    @(15, 'MATIC', 'Polygon', 0.782527, 7257898397, 216356049, -0.89416),
    @(16, 'TON', 'Toncoin', 2.03, 7008074556, 10834302, -0.53139),
    @(17, 'WBTC', 'Wrapped Bitcoin', 43041, 6789970297, 54020330, 0.29294),
    @(18, 'ICP', 'Internet Computer', 12.79, 5846010965, 127112967, -5.4726),
    @(19, 'SHIB', 'Shiba Inu', 0.00000904, 5328123409, 68180466, -0.48268),
    @(20, 'LTC', 'Litecoin', 68.06, 5045735459, 196533260, -0.28399),
    @(21, 'DAI', 'Dai', 0.999953, 4883075320, 137171597, 0.07278),
    @(22, 'BCH', 'Bitcoin Cash', 240.64, 4718223991, 148152928, 0.54599),
    @(23, 'UNI', 'Uniswap', 6.05, 4560214336, 50631185, -1.779),
    @(24, 'LEO', 'LEO Token', 4.03, 3734005183, 3879715, 0.03114),
    @(25, 'ETC', 'Ethereum Classic', 25.45, 3642954624, 192279476, 2.34966),
    @(26, 'ATOM', 'Cosmos Hub', 9.09, 3485137010, 89213967, -0.418),
    @(28, 'XLM', 'Stellar', 0.110071, 3125936788, 47617622, -1.61404),
    @(27, 'APT', 'Aptos', 9.039999999999999, 3056981325, 91573729, -5.98236),
    @(29, 'XMR', 'Monero', 165.81, 3007250631, 59004838, -0.85444),
    @(30, 'INJ', 'Injective', 33.27, 2939236856, 78171856, -0.44431),
    @(31, 'IMX', 'Immutable', 2.14, 2910757206, 56924594, -1.08566),
    @(33, 'OKB', 'OKB', 48.31, 2897575013, 3944643, -1.23228),
    @(32, 'NEAR', 'NEAR Protocol', 2.81, 2895186098, 75562330, -0.8339800000000001),
    @(34, 'OP', 'Optimism', 2.98, 2851474422, 67844632, -1.14699),
    @(35, 'TIA', 'Celestia', 17.08, 2762642855, 97207793, -5.60041),
    @(36, 'TAO', 'Bittensor', 428.69, 2646926542, 6352589, -3.85166),
    @(37, 'FDUSD', 'First Digital USD', 0.99978, 2583312207, 2083425899, 0.01433),
    @(38, 'FIL', 'Filecoin', 5.04, 2540025034, 114450557, -3.83839),
    @(39, 'LDO', 'Lido DAO', 2.75, 2449495979, 36056375, -0.26438),
    @(40, 'HBAR', 'Hedera', 0.068771, 2314048921, 21909276, -0.51732),
    @(41, 'KAS', 'Kaspa', 0.10075, 2272704443, 9324132, -0.80327),
    @(42, 'ARB', 'Arbitrum', 1.77, 2253347963, 169051386, -1.53622),
    @(44, 'CRO', 'Cronos', 0.081201, 2150843073, 4892980, -0.50139),
    @(43, 'STX', 'Stacks', 1.5, 2146889257, 22016309, -2.41355),
    @(45, 'VET', 'VeChain', 0.02826397, 2053224914, 28278693, -2.07937),
    @(46, 'MNT', 'Mantle', 0.593969, 1880414962, 21628765, 0.36955),
    @(49, 'RNDR', 'Render', 4.97, 1869635138, 116083089, 2.42232),
    @(47, 'MKR', 'Maker', 2020.5, 1864293643, 43379023, 0.75641),
    @(48, 'SUI', 'Sui', 1.51, 1759208928, 263661177, 7.42558),
    @(50, 'SEI', 'Sei', 0.619074, 1499436273, 100420024, 1.4926),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
